# Applies the diagnoses.xlsx data update described by the commit:
# "added check for spike value (WIP) and sends alert email to admin"
#
# Row-by-row changes to the "Diagnoses" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("E2").Value = 96.98
$ws.Cells.Item(2, 7).ClearContents()        # G2 "Pains" -> (cleared)
$ws.Range("I2").Value = 1
$ws.Range("L2").Value = 0

# --- Row 3 ---------------------------------------------------------------
$ws.Range("E3").Value = 102.2
$ws.Range("G3").Value = "Chest Pain,Fever,Aches,Sore Throat"
$ws.Range("H3").Value = "Asthma,Dementia"
$ws.Range("I3").Value = 2
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 1
$ws.Range("Q3").Value = "Very High Risk"

# --- Row 4 ---------------------------------------------------------------
$ws.Range("B4").Value = "Jane"
$ws.Range("C4").Value = "Doe"
$ws.Range("E4").Value = 116.6
$ws.Range("F4").Value = 68
$ws.Range("G4").Value = "Difficulty Breathing,Fever,Tiredness,Sore Throat,Conjunctivitis,Headache"
$ws.Range("H4").Value = "Cancer,Asthma,Pulmonary Hypertension,Dementia"
$ws.Range("I4").Value = 4
$ws.Range("J4").Value = 1
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 3
$ws.Range("Q4").Value = "Very High Risk"

# --- Row 5 ---------------------------------------------------------------
$ws.Range("B5").Value = "Cha"
$ws.Range("C5").Value = "O"
$ws.Range("E5").Value = 111.2
$ws.Range("F5").Value = 90
$ws.Range("G5").Value = "Difficulty Breathing,Loss of Movement,Fever"
$ws.Range("H5").Value = "Asthma,Pulmonary Hypertension"
$ws.Range("I5").Value = 2
$ws.Range("K5").Value = 1
